$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Update the "ScriptLatestRunVersion" shared text (column AJ) for all data rows (2-80)
$ws.Range("AJ2:AJ80").Value = "IndicatorQuantiles.R, Git Commit ID: d77a77d64f72a744c78cd38270c72c5d9c8cd498"

# Update the "pid" numeric value (column AH) from 19980 to 11992 for all data rows (2-80)
$ws.Range("AH2:AH80").Value = 11992
